$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header/value text (import fix) ---
$ws.Range("B1").Value = "Asset Name"
$ws.Range("B2").Value = "Motor cycle (Required)"
$ws.Range("K2").Value = "Mantimin (Required)"
$ws.Range("L2").Value = "HRGA (Required)"

# --- Adjust column widths to fit new, longer content ---
$ws.Columns.Item(2).ColumnWidth = 20.5546875
$ws.Columns.Item(11).ColumnWidth = 18.5546875

# --- Update view (scroll position / active selection) ---
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("L7").Select()
